$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text columns (Coin name / Link) - plain strings, not numeric-looking
$textUpdates = @{
    'B8' = 'MXToken'
    'C8' = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
    'B9' = 'WazirX'
    'C9' = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
    'B10' = 'LiechtensteinCryptoassetsExchange'
    'C10' = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
    'B11' = 'MandalaExchangeToken'
    'C11' = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
    'B12' = 'BitrueCoin'
    'C12' = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
    'B13' = 'BitMartToken'
    'C13' = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
    'B14' = 'BitForexToken'
    'C14' = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
    'B15' = 'TigerCash'
    'C15' = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
    'B16' = 'LEO'
    'C16' = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
    'B17' = 'GateToken'
    'C17' = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
}

# Numeric-looking columns (Price / Volume%) stored as literal text -
# prefix with an apostrophe so Excel keeps them as text, matching the
# original inline-string cells instead of converting to numbers/percentages.
$numericTextUpdates = @{
    'E2' = '-0.65%'
    'D3' = '31.44'
    'E3' = '-1.84%'
    'D4' = '5.138'
    'E4' = '-2.11%'
    'D5' = '0.07368'
    'E5' = '-1.25%'
    'D6' = '2.443'
    'E6' = '58.51%'
    'D7' = '7.917'
    'E7' = '0.79%'
    'D8' = '0.9208'
    'E8' = '0.20%'
    'D9' = '0.1737'
    'E9' = '3.26%'
    'D10' = '0.07448'
    'E10' = '-6.60%'
    'D11' = '0.08110'
    'E11' = '1.45%'
    'D12' = '0.03033'
    'E12' = '-0.12%'
    'D13' = '0.09931'
    'E13' = '0.28%'
    'D14' = '0.001512'
    'E14' = '1.10%'
    'D15' = '0.006073'
    'E15' = '-2.26%'
    'D16' = '3.452'
    'E16' = '-0.55%'
    'D17' = '3.755'
    'E17' = '-1.09%'
    'E18' = '-0.09%'
    'D19' = '0.3291'
    'E19' = '-1.05%'
    'D20' = '0.1338'
    'E20' = '-0.06%'
    'D21' = '4.656'
    'E21' = '3.75%'
    'D22' = '0.04653'
    'E22' = '1.10%'
    'D24' = '0.001226'
    'E24' = '0.70%'
    'D25' = '0.004480'
    'E25' = '0.71%'
    'D26' = '0.0001301'
    'E26' = '-7.00%'
    'E27' = '7.13%'
    'D39' = '0.01729'
    'E39' = '-1.40%'
    'D40' = '0.04520'
    'E40' = '0.57%'
    'D41' = '0.007184'
    'E41' = '0.26%'
    'D42' = '0.1344'
    'E42' = '-0.45%'
    'D43' = '0.002222'
    'E43' = '0.60%'
    'D44' = '0.01093'
    'D45' = '0.00006273'
    'E45' = '1.88%'
    'D46' = '0.01000'
    'E46' = '-23.11%'
    'E47' = '171.86%'
}

foreach ($cell in $textUpdates.Keys) {
    $ws.Range($cell).Value = $textUpdates[$cell]
}

foreach ($cell in $numericTextUpdates.Keys) {
    $ws.Range($cell).Value = "'" + $numericTextUpdates[$cell]
}
